$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the old row 147, pushing the
# remaining rows (old 147-152) down by one (new 148-153).
$ws.Rows.Item(147).Insert()

$ws.Cells.Item(147, 1).Value = 10
$ws.Cells.Item(147, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(147, 3).Value = "La Araucanía"
$ws.Cells.Item(147, 4).Value = 44509
$ws.Cells.Item(147, 5).Value = 9
$ws.Cells.Item(147, 6).Value = 100112005
$ws.Cells.Item(147, 7).Value = "Puerro"
$ws.Cells.Item(147, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 40
$ws.Cells.Item(147, 11).Value = 7000
$ws.Cells.Item(147, 12).Value = 7000
$ws.Cells.Item(147, 13).Value = 7000
$ws.Cells.Item(147, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(147, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(147, 16).Value = 583
$ws.Cells.Item(147, 17).Value = 12
$ws.Cells.Item(147, 18).Value = "Hortaliza"
